# "case b - heatsink script"
# Insert a new row above the old row 5 (Final Cell Temperature block),
# turn the new row into a "Current" row with its own formula, and trim
# the block that got pushed down to row 6 so only the G-column formula
# (now referencing the shifted J20 cell) survives.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 5 - shifts old rows 5:25 down to 6:26
#    and auto-adjusts every formula reference (K11->K12, J19->J20, etc.)
$ws.Rows.Item(5).EntireRow.Insert()

# 2. New row 5 holds a "Current" label + a new current-draw formula
$ws.Range("A5").Value = "Current"
$ws.Range("B5").Formula = "=B4/504*1000"

# 3. Give the J5:M5 placeholder cells the same green fill used elsewhere
#    in this block (matches the fill used by fillId 3 / style index 3).
$ws.Range("J5:M5").Interior.Color = 5296274

# 4. Row 6 (the old row 5, shifted down) keeps its "Final Cell
#    Temperature" label in A6 and the peak-power formula in G6, but the
#    per-voltage-level (B:F) ohmic-heating formulas it used to carry are
#    no longer needed there, so clear them out.
$ws.Range("B6:F6").ClearContents()

# 5. Restore the selection to where the author left off editing.
$ws.Range("A6").Select()
